# Adds the new match row (Istra 1961 vs Slaven Belupo) at the bottom of the
# croatia_hnl_2023-2024 sheet, as produced by the nightly scraping script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow = 66
$newRow = 67

# Copy formatting only (number formats, font, border, alignment) from the
# last existing data row so the new row matches the established style
# (bold/bordered index column, datetime-formatted match-date column, etc.)
$ws.Range("A" + $sourceRow + ":V" + $sourceRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's values
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = "croatia"
$ws.Range("C67").Value = "hnl"
$ws.Range("D67").Value = "2023-2024"
$ws.Range("E67").Value = 45235.71527777778
$ws.Range("F67").Value = "Istra 1961"
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = "Slaven Belupo"
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1.98
$ws.Range("K67").Value = "29/10/2023 17:13"
$ws.Range("L67").Value = 2.23
$ws.Range("M67").Value = "05/11/2023 17:04"
$ws.Range("N67").Value = 3.34
$ws.Range("O67").Value = "29/10/2023 17:13"
$ws.Range("P67").Value = 3.1
$ws.Range("Q67").Value = "05/11/2023 16:59"
$ws.Range("R67").Value = 3.73
$ws.Range("S67").Value = "29/10/2023 17:13"
$ws.Range("T67").Value = 3.66
$ws.Range("U67").Value = "05/11/2023 17:04"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/croatia/hnl/istra-1961-slaven-belupo/j1aQWaZa/"
